$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2250.7
$ws.Range("I88").Value = 1867.6666
$ws.Range("J88").Value = 2414.8572
$ws.Range("K88").Value = 1867.6666
$ws.Range("L88").Value = 2414.8572
$ws.Range("M88").Value = -1461.6666
$ws.Range("N88").Value = -3226.8572
$ws.Range("H91").Value = 2250.7
$ws.Range("I91").Value = 1867.6666
$ws.Range("J91").Value = 2414.8572
$ws.Range("K91").Value = 1867.6666
$ws.Range("L91").Value = 2414.8572
$ws.Range("M91").Value = -463.6666
$ws.Range("N91").Value = -5222.8572
$ws.Range("H113").Value = 4318.3335
$ws.Range("I113").Value = 2755
$ws.Range("J113").Value = 5100
$ws.Range("K113").Value = 2755
$ws.Range("L113").Value = 5100
$ws.Range("M113").Value = 499
$ws.Range("N113").Value = -11608
$ws.Range("H129").Value = 719.5806
$ws.Range("I129").Value = 297.4
$ws.Range("J129").Value = 920.619
$ws.Range("K129").Value = 892.1999999999999
$ws.Range("L129").Value = 2761.857
$ws.Range("M129").Value = 4107.8
$ws.Range("N129").Value = -12761.857
$ws.Range("H135").Value = 3102.6072
$ws.Range("I135").Value = 213.66667
$ws.Range("J135").Value = 8302.700000000001
$ws.Range("K135").Value = 1923.00003
$ws.Range("L135").Value = 74724.3
$ws.Range("M135").Value = 611.9999699999998
$ws.Range("N135").Value = -79794.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3157
$ws.Range("I63").Value = 3552.5
$ws.Range("J63").Value = 1575
$ws.Range("K63").Value = 3552.5
$ws.Range("L63").Value = 1575
$ws.Range("M63").Value = -2866.5
$ws.Range("N63").Value = -2947
$ws.Range("H66").Value = 3157
$ws.Range("I66").Value = 3552.5
$ws.Range("J66").Value = 1575
$ws.Range("K66").Value = 17762.5
$ws.Range("L66").Value = 7875
$ws.Range("M66").Value = -14330.5
$ws.Range("N66").Value = -14739
$ws.Range("H74").Value = 813.61536
$ws.Range("I74").Value = 774.3333
$ws.Range("K74").Value = 774.3333
$ws.Range("M74").Value = 99.66669999999999
$ws.Range("H77").Value = 813.61536
$ws.Range("I77").Value = 774.3333
$ws.Range("K77").Value = 3871.6665
$ws.Range("M77").Value = 496.3334999999997
$ws.Range("H88").Value = 2584.7144
$ws.Range("I88").Value = 1426.5
$ws.Range("J88").Value = 3048
$ws.Range("K88").Value = 1426.5
$ws.Range("L88").Value = 3048
$ws.Range("M88").Value = -1020.5
$ws.Range("N88").Value = -3860
$ws.Range("H91").Value = 2584.7144
$ws.Range("I91").Value = 1426.5
$ws.Range("J91").Value = 3048
$ws.Range("K91").Value = 1426.5
$ws.Range("L91").Value = 3048
$ws.Range("M91").Value = -22.5
$ws.Range("N91").Value = -5856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 487.5909
$ws.Range("I64").Value = 474
$ws.Range("J64").Value = 489.73685
$ws.Range("K64").Value = 474
$ws.Range("L64").Value = 489.73685
$ws.Range("M64").Value = -249
$ws.Range("N64").Value = -939.73685
$ws.Range("H67").Value = 487.5909
$ws.Range("I67").Value = 474
$ws.Range("J67").Value = 489.73685
$ws.Range("K67").Value = 474
$ws.Range("L67").Value = 489.73685
$ws.Range("M67").Value = 306
$ws.Range("N67").Value = -2049.73685
$ws.Range("H86").Value = 3324
$ws.Range("I86").Value = 2692.05
$ws.Range("J86").Value = 5129.5713
$ws.Range("K86").Value = 2692.05
$ws.Range("L86").Value = 5129.5713
$ws.Range("M86").Value = -1569.05
$ws.Range("N86").Value = -7375.5713
$ws.Range("H89").Value = 3324
$ws.Range("I89").Value = 2692.05
$ws.Range("J89").Value = 5129.5713
$ws.Range("K89").Value = 13460.25
$ws.Range("L89").Value = 25647.8565
$ws.Range("M89").Value = -7844.25
$ws.Range("N89").Value = -36879.85649999999
$ws.Range("H105").Value = 1773.8572
$ws.Range("I105").Value = 1719.3529
$ws.Range("J105").Value = 2005.5
$ws.Range("K105").Value = 1719.3529
$ws.Range("L105").Value = 2005.5
$ws.Range("M105").Value = 27.64709999999991
$ws.Range("N105").Value = -5499.5
$ws.Range("H107").Value = 1912.6428
$ws.Range("I107").Value = 1553.8948
$ws.Range("J107").Value = 2670
$ws.Range("K107").Value = 1553.8948
$ws.Range("L107").Value = 2670
$ws.Range("M107").Value = 366.1052
$ws.Range("N107").Value = -6510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1768.4474
$ws.Range("I58").Value = 1764.0333
$ws.Range("J58").Value = 1785
$ws.Range("K58").Value = 1764.0333
$ws.Range("L58").Value = 1785
$ws.Range("M58").Value = -1561.0333
$ws.Range("N58").Value = -2191
$ws.Range("H94").Value = 161329.4
$ws.Range("I94").Value = 250916.38
$ws.Range("J94").Value = 119170.82
$ws.Range("K94").Value = 250916.38
$ws.Range("L94").Value = 119170.82
$ws.Range("M94").Value = -250465.38
$ws.Range("N94").Value = -120072.82
$ws.Range("H132").Value = 2014.4348
$ws.Range("I132").Value = 1575.4736
$ws.Range("J132").Value = 4099.5
$ws.Range("K132").Value = 4726.4208
$ws.Range("L132").Value = 12298.5
$ws.Range("M132").Value = -2196.4208
$ws.Range("N132").Value = -17358.5
$ws.Range("H136").Value = 1768.4474
$ws.Range("I136").Value = 1764.0333
$ws.Range("J136").Value = 1785
$ws.Range("K136").Value = 5292.0999
$ws.Range("L136").Value = 5355
$ws.Range("M136").Value = -2742.0999
$ws.Range("N136").Value = -10455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 197
$ws.Range("I2").Value = 11.4
$ws.Range("J2").Value = 727.2857
$ws.Range("K2").Value = 68.40000000000001
$ws.Range("L2").Value = 4363.7142
$ws.Range("M2").Value = 44.59999999999999
$ws.Range("N2").Value = -4589.7142
$ws.Range("H5").Value = 22564.787
$ws.Range("I5").Value = 471.6842
$ws.Range("J5").Value = 37556.535
$ws.Range("K5").Value = 1415.0526
$ws.Range("L5").Value = 112669.605
$ws.Range("M5").Value = -1303.0526
$ws.Range("N5").Value = -112893.605
$ws.Range("H23").Value = 97.5
$ws.Range("I23").Value = 30.166666
$ws.Range("J23").Value = 137.9
$ws.Range("K23").Value = 90.49999800000001
$ws.Range("L23").Value = 413.7
$ws.Range("M23").Value = 144.500002
$ws.Range("N23").Value = -883.7
$ws.Range("H34").Value = 1330.7
$ws.Range("I34").Value = 780
$ws.Range("J34").Value = 1627.2307
$ws.Range("K34").Value = 2340
$ws.Range("L34").Value = 4881.6921
$ws.Range("M34").Value = -2256
$ws.Range("N34").Value = -5049.6921
$ws.Range("H39").Value = 1123.4
$ws.Range("I39").Value = 1195
$ws.Range("J39").Value = 1115.4445
$ws.Range("K39").Value = 3585
$ws.Range("L39").Value = 3346.3335
$ws.Range("M39").Value = -3291
$ws.Range("N39").Value = -3934.3335
$ws.Range("H55").Value = 2694.6667
$ws.Range("I55").Value = 2704
$ws.Range("J55").Value = 2694.1177
$ws.Range("K55").Value = 8112
$ws.Range("L55").Value = 8082.353099999999
$ws.Range("M55").Value = -7935
$ws.Range("N55").Value = -8436.3531
$ws.Range("H135").Value = 22564.787
$ws.Range("I135").Value = 471.6842
$ws.Range("J135").Value = 37556.535
$ws.Range("K135").Value = 4245.1578
$ws.Range("L135").Value = 338008.8150000001
$ws.Range("M135").Value = -1710.1578
$ws.Range("N135").Value = -343078.8150000001
$ws.Range("H137").Value = 5379776.5
$ws.Range("I137").Value = 2259
$ws.Range("J137").Value = 6948219
$ws.Range("K137").Value = 6777
$ws.Range("L137").Value = 20844657
$ws.Range("M137").Value = -1677
$ws.Range("N137").Value = -20854857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3821.3713
$ws.Range("I70").Value = 3793.6
$ws.Range("J70").Value = 3890.8
$ws.Range("K70").Value = 3793.6
$ws.Range("L70").Value = 3890.8
$ws.Range("M70").Value = -3523.6
$ws.Range("N70").Value = -4430.8
$ws.Range("H73").Value = 3821.3713
$ws.Range("I73").Value = 3793.6
$ws.Range("J73").Value = 3890.8
$ws.Range("K73").Value = 3793.6
$ws.Range("L73").Value = 3890.8
$ws.Range("M73").Value = -2857.6
$ws.Range("N73").Value = -5762.8
$ws.Range("H80").Value = 2405.8823
$ws.Range("I80").Value = 2366.6667
$ws.Range("J80").Value = 2414.2856
$ws.Range("K80").Value = 2366.6667
$ws.Range("L80").Value = 2414.2856
$ws.Range("M80").Value = -1368.6667
$ws.Range("N80").Value = -4410.2856
$ws.Range("H83").Value = 2405.8823
$ws.Range("I83").Value = 2366.6667
$ws.Range("J83").Value = 2414.2856
$ws.Range("K83").Value = 11833.3335
$ws.Range("L83").Value = 12071.428
$ws.Range("M83").Value = -6841.333500000001
$ws.Range("N83").Value = -22055.428
$ws.Range("H97").Value = 1287.4615
$ws.Range("I97").Value = 1144.75
$ws.Range("K97").Value = 1144.75
$ws.Range("M97").Value = -648.75
$ws.Range("H126").Value = 12666.667
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 734.1795
$ws.Range("I136").Value = 669.62067
$ws.Range("J136").Value = 921.4
$ws.Range("K136").Value = 2008.86201
$ws.Range("L136").Value = 2764.2
$ws.Range("M136").Value = 541.1379899999999
$ws.Range("N136").Value = -7864.2
